$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Settings string shown in N1 (shared string): -platformcnt -10 -> -1 ---
$ws.Range("N1").Value = "Settings: '-spcnt 0 -platformcnt -1 -personcnt 1'"

# --- 2. Refreshed benchmark data: Creation Time (Awa) / Search Time (Awa), rows 2-31 ---
$ws.Range("K2").Value = 13.044903
$ws.Range("L2").Value = 223.06951100000001
$ws.Range("K3").Value = 12.740627999999999
$ws.Range("L3").Value = 230.943242
$ws.Range("K4").Value = 10.611908
$ws.Range("L4").Value = 217.69246899999999
$ws.Range("K5").Value = 13.360951
$ws.Range("L5").Value = 228.125978
$ws.Range("K6").Value = 14.165108
$ws.Range("L6").Value = 263.494978
$ws.Range("K7").Value = 15.806323000000001
$ws.Range("L7").Value = 240.38362900000001
$ws.Range("K8").Value = 12.866201999999999
$ws.Range("L8").Value = 263.86596400000002
$ws.Range("K9").Value = 12.800094
$ws.Range("L9").Value = 241.20740599999999
$ws.Range("K10").Value = 22.639842000000002
$ws.Range("L10").Value = 254.00538800000001
$ws.Range("K11").Value = 18.587063000000001
$ws.Range("L11").Value = 259.46514000000002
$ws.Range("K12").Value = 15.831378000000001
$ws.Range("L12").Value = 256.45106199999998
$ws.Range("K13").Value = 24.456137999999999
$ws.Range("L13").Value = 373.45311600000002
$ws.Range("K14").Value = 18.024093000000001
$ws.Range("L14").Value = 298.67139100000003
$ws.Range("K15").Value = 22.389600999999999
$ws.Range("L15").Value = 498.03247299999998
$ws.Range("K16").Value = 36.025244999999998
$ws.Range("L16").Value = 363.73079100000001
$ws.Range("K17").Value = 39.111165
$ws.Range("L17").Value = 335.15546399999999
$ws.Range("K18").Value = 22.942609999999998
$ws.Range("L18").Value = 292.869933
$ws.Range("K19").Value = 37.693326999999996
$ws.Range("L19").Value = 288.34444000000002
$ws.Range("K20").Value = 57.586238999999999
$ws.Range("L20").Value = 716.17199800000003
$ws.Range("K21").Value = 74.960425999999998
$ws.Range("L21").Value = 333.04666600000002
$ws.Range("K22").Value = 263.34917899999999
$ws.Range("L22").Value = 1439.6864519999999
$ws.Range("K23").Value = 102.681488
$ws.Range("L23").Value = 2397.771002
$ws.Range("K24").Value = 214.71129500000001
$ws.Range("L24").Value = 1876.7702609999999
$ws.Range("K25").Value = 370.14774299999999
$ws.Range("L25").Value = 3361.3008140000002
$ws.Range("K26").Value = 51.258637999999998
$ws.Range("L26").Value = 718.56665699999996
$ws.Range("K27").Value = 402.48123399999997
$ws.Range("L27").Value = 2124.5168720000001
$ws.Range("K28").Value = 2347.0795589999998
$ws.Range("L28").Value = 4635.1813030000003
$ws.Range("K29").Value = 74.959519999999998
$ws.Range("L29").Value = 5104.0190309999998
$ws.Range("K30").Value = 144.595731
$ws.Range("L30").Value = 7838.3213779999996
$ws.Range("K31").Value = 62.140711000000003
$ws.Range("L31").Value = 9029.3701259999998

# --- 3. New data rows 30 and 31 gain H:M columns (two more x=290/300 samples) ---
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 290
$ws.Range("J30").Value = 290
$ws.Range("M30").Value = 290
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 300
$ws.Range("J31").Value = 300
$ws.Range("M31").Value = 300

# --- 4. New helper cell R3: concatenated two-line chart caption ---
$ws.Range("R3").Formula = "=G1&CHAR(10)&N1"

# --- 5. Re-point the chart title at the concatenated caption in R3 ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$title = $chart.ChartTitle
$title.Text = $ws.Range("R3").Value2

# --- 6. Selection moves to R4, matching the saved view state ---
$ws.Range("R4").Select()
